$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column H with header and dataset id values
$ws.Range("H1").Value = "Dataset id"
$ws.Range("H2").Value = 725
$ws.Range("H3").Value = 725
$ws.Range("H4").Value = 726
$ws.Range("H5").Value = 727
$ws.Range("H6").Value = 729
$ws.Range("H7").Value = 728

# Move the active selection to H8, matching the post-edit state
$ws.Range("H8").Select()
